# "j'ai oublie de enregistre le excel"
# Finish the audit updates that were made in the live session but never saved:
#  - mark several rows as done ("FAIT") in column H
#  - flip the "texte en image" accessibility row back to "non corrige" (FALSE)
#  - add the two new SEO rows (texte en image / texture inutile) that were
#    filled in at the bottom of the SEO table
#  - leave the cursor where the user last left it (H20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark rows as done ---------------------------------------------------
$ws.Range("H2").Value = "FAIT"
$ws.Range("H3").Value = "FAIT"
$ws.Range("H6").Value = "FAIT"
$ws.Range("H9").Value = "FAIT"
$ws.Range("H18").Value = "FAIT"

# --- Accessibilité: "texte en image" row wasn't actually fixed yet -------
$ws.Range("E6").Value = $false

# --- New SEO findings added at the bottom of the table -------------------
$ws.Range("A26").Value = "SEO"
$ws.Range("B26").Value = "texte en image"
$ws.Range("C26").Value = "mot clé dans image non`nexploité"
$ws.Range("D26").Value = "remplacer image par du `ntexte"
$ws.Range("E26").Value = $true
$ws.Range("H26").Value = "FAIT"

$ws.Range("A27").Value = "SEO"
$ws.Range("B27").Value = "texture inutile page 2"
$ws.Range("C27").Value = "ralentit chargement de la `npage et mauvais pour le `ncontraste"
$ws.Range("D27").Value = "eviter de surcharger en `ntexture pour faire un design non pertinent et penser au`ncontraste"
$ws.Range("E27").Value = $true
$ws.Range("H27").Value = "FAIT"

$ws.Range("A26").VerticalAlignment = -4160
$ws.Range("B26").VerticalAlignment = -4160
$ws.Range("E26").VerticalAlignment = -4160
$ws.Range("A27").VerticalAlignment = -4160
$ws.Range("E27").VerticalAlignment = -4160

$ws.Range("C26").HorizontalAlignment = -4131
$ws.Range("C26").VerticalAlignment = -4160
$ws.Range("C26").WrapText = $true

$ws.Range("D26").VerticalAlignment = -4160
$ws.Range("D26").WrapText = $true

$ws.Range("B27").VerticalAlignment = -4160
$ws.Range("B27").WrapText = $true

$ws.Range("C27").HorizontalAlignment = -4131
$ws.Range("C27").VerticalAlignment = -4160
$ws.Range("C27").WrapText = $true

$ws.Range("D27").VerticalAlignment = -4160
$ws.Range("D27").WrapText = $true

$ws.Rows("26:26").RowHeight = 33
$ws.Rows("27:27").RowHeight = 82.5

# --- Restore the cursor position left by the user -------------------------
$ws.Range("H20").Select()
